# Applies the commit's data update to the "展览" and "全部类型" worksheets.
# The change removes the two oldest events (rows 2-3, the TCD events) which
# shifts every remaining event up by two rows, and also refreshes the
# "想去人数" (F column) view-counts that ticked up between scrapes.

$wb = $excel.ActiveWorkbook

# New F-column ("想去人数") values for the surviving rows, in final row order
# (worksheet rows 2..26 after the two old rows have been removed).
$newF = @(7,1064,125,80,50,52,11054,4253,23,18,13,2487,1063,73,11,150,473,11194,11021,13,35,11,7,30,16)

$targetSheets = @("展览", "全部类型")

foreach ($sheetName in $targetSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Remove the two obsolete rows (old row 2 and old row 3); Excel shifts
    # everything below up automatically.
    $ws.Rows.Item(2).Delete()
    $ws.Rows.Item(2).Delete()

    # Re-sequence the serial-number column (A) to 1..25 for the 25 data rows.
    for ($i = 0; $i -lt 25; $i++) {
        $r = $i + 2
        $ws.Cells.Item($r, 1).Value2 = $i + 1
    }

    # Refresh the "want to go" counts that changed since the previous scrape.
    for ($i = 0; $i -lt $newF.Length; $i++) {
        $r = $i + 2
        $ws.Cells.Item($r, 6).Value2 = $newF[$i]
    }
}
